$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for price cells whose new values would otherwise be
# auto-coerced into numbers by Excel (e.g. "1.00" -> 1, "103.00" -> 103).
$textCells = @("D4", "D5", "D6", "D7", "D10", "D15", "D17", "D19", "D20", "D22", "D23", "D26", "D28", "D29", "D32", "D34", "D35", "D36", "D37", "D38", "D41", "D44", "D45", "D48", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated coin data (price + 1h volume change), including the
# three pairs of rows whose coin/link/price/volume content was reordered.
$ws.Range("D2").Value = "43.024.10"
$ws.Range("E2").Value = "  +0.92%  "
$ws.Range("D3").Value = "2.535.52"
$ws.Range("E3").Value = "  +0.72%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "318.28"
$ws.Range("E5").Value = "  +4.56%  "
$ws.Range("D6").Value = "96.08"
$ws.Range("E6").Value = "  -0.72%  "
$ws.Range("D7").Value = "0.580"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("E9").Value = "  -0.36%  "
$ws.Range("D10").Value = "36.41"
$ws.Range("E10").Value = "  -0.43%  "
$ws.Range("E11").Value = "  +0.24%  "
$ws.Range("E12").Value = "  -0.36%  "
$ws.Range("E13").Value = "  -0.10%  "
$ws.Range("D14").Value = "2.926.62"
$ws.Range("E14").Value = "  +0.75%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "15.40"
$ws.Range("E15").Value = "  +2.64%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.515.61"
$ws.Range("E16").Value = "  -1.62%  "
$ws.Range("D17").Value = "0.855"
$ws.Range("E17").Value = "  -0.80%  "
$ws.Range("D18").Value = "43.084.27"
$ws.Range("E18").Value = "  +0.90%  "
$ws.Range("D19").Value = "13.10"
$ws.Range("E19").Value = "  +1.77%  "
$ws.Range("D20").Value = "6.67"
$ws.Range("E20").Value = "  +3.59%  "
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("D22").Value = "70.28"
$ws.Range("E22").Value = "  -1.21%  "
$ws.Range("D23").Value = "252.87"
$ws.Range("E23").Value = "  +0.53%  "
$ws.Range("E24").Value = "  +2.10%  "
$ws.Range("E25").Value = "  -0.25%  "
$ws.Range("D26").Value = "27.10"
$ws.Range("E26").Value = "  +0.45%  "
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("D28").Value = "2.44"
$ws.Range("E28").Value = "  +5.07%  "
$ws.Range("D29").Value = "40.02"
$ws.Range("E29").Value = "  +5.04%  "
$ws.Range("E30").Value = "  -0.23%  "
$ws.Range("E31").Value = "  +2.00%  "
$ws.Range("D32").Value = "154.92"
$ws.Range("E32").Value = "  -1.38%  "
$ws.Range("E33").Value = "  +2.98%  "
$ws.Range("D34").Value = "3.32"
$ws.Range("E34").Value = "  +0.66%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "0.0795"
$ws.Range("E35").Value = "  +0.78%  "
$ws.Range("B36").Value = "Celestia"
$ws.Range("C36").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D36").Value = "18.90"
$ws.Range("E36").Value = "  +2.28%  "
$ws.Range("D37").Value = "2.60"
$ws.Range("E37").Value = "  -0.51%  "
$ws.Range("D38").Value = "0.113"
$ws.Range("E38").Value = "  -2.85%  "
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("E40").Value = "  -1.27%  "
$ws.Range("D41").Value = "2.32"
$ws.Range("E41").Value = "  +12.89%  "
$ws.Range("E42").Value = "  +2.19%  "
$ws.Range("E43").Value = "  -0.14%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.43%  "
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").Value = "3.32"
$ws.Range("E45").Value = "  -2.27%  "
$ws.Range("D46").Value = "2.025.19"
$ws.Range("E47").Value = "  +0.99%  "
$ws.Range("D48").Value = "8.86"
$ws.Range("E48").Value = "  -1.40%  "
$ws.Range("D49").Value = "2.780.50"
$ws.Range("E49").Value = "  +0.66%  "
$ws.Range("D50").Value = "74.04"
$ws.Range("E50").Value = "  +2.66%  "
$ws.Range("D51").Value = "103.00"
$ws.Range("E51").Value = "  +1.61%  "
